$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Formula = "=SUM(Sheet1:Sheet3!A1:B2)"
$ws.Range("I3").Formula = "=AVERAGE(Sheet1:Sheet3!A1:B2)"
$ws.Range("I4").Formula = "=MIN(Sheet1:Sheet3!A`$1:B`$2)"
$ws.Range("I5").Formula = "=MAX(Sheet1:Sheet3!A`$1:B`$2)"
$ws.Range("I6").Formula = "=COUNT(Sheet1:Sheet3!`$A`$1:`$B`$2)"

$ws.Range("H3").Select()
